# -----------------------------------------------------------------------
# "added generic algorithm diagrams"
#
# Renames the blank "Sheet1" to "Generic Layer Test", adds a new "L0"
# defined name pointing at its input cell, and fills that sheet with a
# small worked example of forward- and back-propagation through a
# generic single-input "layer" (Z = L*W, L = sigmoid(Z), an error /
# back-prop block, and an "Initial" / "pass 1" convergence table).
# Also updates column widths and the active sheet/cell selections to
# match the saved view state.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- rename the blank "Sheet1" -> "Generic Layer Test" ---
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Name = "Generic Layer Test"

# --- new workbook-level name: the "layer 0" input cell ---
$wb.Names.Add("L0", "='Generic Layer Test'!`$E`$7")

# --- column widths for the new table (C:F) ---
$ws.Columns.Item(3).ColumnWidth = 11.85546875
$ws.Columns.Item(4).ColumnWidth = 11.42578125
$ws.Columns.Item(5).ColumnWidth = 14.28515625
$ws.Columns.Item(6).ColumnWidth = 10.85546875

# --- header / label text, written first and in first-use reading order so
#     the shared-string table lands in the same order as the authored file ---
$ws.Range("B6").Value = "Layer"
$ws.Range("C6").Value = "W"
$ws.Range("B1").Value = "Input X"
$ws.Range("B2").Value = "Target Y"
$ws.Range("D1").Value = " "
$ws.Range("D6").Value = "Z = L * W"
$ws.Range("A6").Value = "Forward"
$ws.Range("A13").Value = "Back"
$ws.Range("C12").Value = "E"
$ws.Range("E6").Value = "L = A(Z)"
$ws.Range("D12").Value = "C(W)"
$ws.Range("B3").Value = "Rate"
$ws.Range("C18").Value = "Initial"
$ws.Range("D18").Value = "pass 1"
$ws.Range("E12").Value = "dW = E * L"
$ws.Range("F12").Value = "W -= lr * dW"

# --- numeric inputs / formulas for the worked example ---
$ws.Range("C1").Value = 4
$ws.Range("C1").Style = "Input"
$ws.Range("C2").Value = 7
$ws.Range("C2").Style = "Input"
$ws.Range("C3").Value = 0.01
$ws.Range("C3").Style = "Input"
$ws.Range("B7").Value = 0
$ws.Range("E7").Formula = "=C1"
$ws.Range("E7").Style = "Calculation"
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 0.51721518182005644
$ws.Range("D8").Formula = "=L0*C8"
$ws.Range("D8").Style = "Calculation"
$ws.Range("E8").Formula = "=1/(1+EXP(-D8))"
$ws.Range("E8").Style = "Calculation"
$ws.Range("B9").Value = 2
$ws.Range("C9").Value = 0.50109152103123156
$ws.Range("D9").Formula = "=E8*C9"
$ws.Range("D9").Style = "Calculation"
$ws.Range("E9").Formula = "=1/(1+EXP(-D9))"
$ws.Range("E9").Style = "Calculation"
$ws.Range("B10").Value = 3
$ws.Range("C10").Value = 0.51720999828656755
$ws.Range("D10").Formula = "=E9*C10"
$ws.Range("D10").Style = "Calculation"
$ws.Range("E10").Formula = "=1/(1+EXP(-D10))"
$ws.Range("E10").Style = "Calculation"
$ws.Range("B13").Value = 3
$ws.Range("C13").Formula = "=(E10-`$C`$2)*(1/(1+EXP(-E10))*(1-(1/(1+EXP(-E10)))))"
$ws.Range("D13").Formula = "=C13*C9"
$ws.Range("E13").Formula = "=C13*E10"
$ws.Range("F13").Formula = "=C10-`$C`$3*E13"
$ws.Range("B14").Value = 2
$ws.Range("C14").Formula = "=(C13*C10)*(1/(1+EXP(-D9))*(1-(1/(1+EXP(-D9)))))"
$ws.Range("D14").Formula = "=C14*E8"
$ws.Range("E14").Formula = "=C14*E9"
$ws.Range("F14").Formula = "=C19-`$C`$3*E14"
$ws.Range("B15").Value = 1
$ws.Range("C15").Formula = "=(C14*C9)*(1/(1+EXP(-D8))*(1-(1/(1+EXP(-D8)))))"
$ws.Range("D15").Formula = "=C15*L0"
$ws.Range("E15").Formula = "=C15*E8"
$ws.Range("F15").Formula = "=C8-`$C`$3*E15"
$ws.Range("B16").Value = 0
$ws.Range("E18").Value = 2
$ws.Range("F18").Value = 3
$ws.Range("C19").Value = 0.5
$ws.Range("C19").Style = "Input"
$ws.Range("D19").Value = 0.50851784753127982
$ws.Range("E19").Value = 0.50860039221460152
$ws.Range("F19").Value = 0.51713063151095084
$ws.Range("G19").Value = 0.51721518182005644
$ws.Range("C20").Value = 0.5
$ws.Range("C20").Style = "Input"
$ws.Range("D20").Value = 0.50107272564531102
$ws.Range("E20").Value = 0.50107336469369013
$ws.Range("F20").Value = 0.50109104405172933
$ws.Range("G20").Value = 0.50109152103123156
$ws.Range("C21").Value = 0.5
$ws.Range("C21").Style = "Input"
$ws.Range("D21").Value = 0.50008153441100756
$ws.Range("E21").Value = 0.50859773614299053
$ws.Range("F21").Value = 0.50868157871043673
$ws.Range("G21").Value = 0.51720999828656755

# --- selection / active-sheet view state ---
$wsFirst = $wb.Worksheets.Item("activation test data (2)")
[void]$wsFirst.Range("F46").Select()

[void]$ws.Activate()
[void]$ws.Range("I10").Select()
